# Actualización automática del mapa: agrega el nuevo registro (fila 46)
# al final de la hoja, replicando el formato de texto de las filas previas
# (los valores "numéricos" de estas columnas se guardan como texto).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 46

# Columnas que deben quedar como texto aunque su contenido sea numérico
# (Caso, F. De Reclamo, Comuna, OT) -> se formatean como texto antes de
# escribir el valor para que no se conviertan en números.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 5).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value  = "-546"
$ws.Cells.Item($row, 2).Value  = "8/5/2025"
$ws.Cells.Item($row, 3).Value  = "Albarellos 3031"
$ws.Cells.Item($row, 4).Value  = "12"
$ws.Cells.Item($row, 5).Value  = "808720857"
$ws.Cells.Item($row, 6).Value  = "INCO"
$ws.Cells.Item($row, 7).Value  = "Pendiente"
$ws.Cells.Item($row, 8).Value  = "Picada"
$ws.Cells.Item($row, 9).Value  = 1
$ws.Cells.Item($row, 10).Value = "Cambio"
$ws.Cells.Item($row, 11).Value = "Sin equipos"
$ws.Cells.Item($row, 12).Value = "Pasante"
$ws.Cells.Item($row, 13).Value = -58.511732
$ws.Cells.Item($row, 14).Value = -34.578688
$ws.Cells.Item($row, 15).Value = "Paternal"
$ws.Cells.Item($row, 16).Value = "Capital Norte"
